$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: Status changed from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both files / both locale columns.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Columns E/F widen to fit the new (longer) status text.
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: record the handback target/file + datetime for both rows.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Row 2 (637a9a6d...)
$zhcn.Range("I2").Value = "637a9a6d-8ca1-4c7f-8b28-254efe9721e2.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c24e5a413eb6017ef38dccc6c1f8676bc1cbad65/e2e/637a9a6d-8ca1-4c7f-8b28-254efe9721e2.md", "", "", "637a9a6d-8ca1-4c7f-8b28-254efe9721e2.md")
$zhcn.Range("J2").Value = "637a9a6d-8ca1-4c7f-8b28-254efe9721e2.ab14a363f11a8f1811c2b2cc7e9f9a0e80da3074.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 18:55:13"

# Row 3 (d74854b4...)
$zhcn.Range("I3").Value = "d74854b4-26bb-4c7f-a523-8082fbbe6f40.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c24e5a413eb6017ef38dccc6c1f8676bc1cbad65/e2e/d74854b4-26bb-4c7f-a523-8082fbbe6f40.md", "", "", "d74854b4-26bb-4c7f-a523-8082fbbe6f40.md")
$zhcn.Range("J3").Value = "d74854b4-26bb-4c7f-a523-8082fbbe6f40.0282d032c43f4d31664d01bc9026ed75a88ebd34.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-31 18:55:13"

# Widen columns to match the new content.
$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# 3. de-de sheet: record the handback target/file + datetime for both rows.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Row 2 (637a9a6d...)
$dede.Range("I2").Value = "637a9a6d-8ca1-4c7f-8b28-254efe9721e2.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c24e5a413eb6017ef38dccc6c1f8676bc1cbad65/e2e/637a9a6d-8ca1-4c7f-8b28-254efe9721e2.md", "", "", "637a9a6d-8ca1-4c7f-8b28-254efe9721e2.md")
$dede.Range("J2").Value = "637a9a6d-8ca1-4c7f-8b28-254efe9721e2.ab14a363f11a8f1811c2b2cc7e9f9a0e80da3074.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 18:55:25"

# Row 3 (d74854b4...)
$dede.Range("I3").Value = "d74854b4-26bb-4c7f-a523-8082fbbe6f40.md"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c24e5a413eb6017ef38dccc6c1f8676bc1cbad65/e2e/d74854b4-26bb-4c7f-a523-8082fbbe6f40.md", "", "", "d74854b4-26bb-4c7f-a523-8082fbbe6f40.md")
$dede.Range("J3").Value = "d74854b4-26bb-4c7f-a523-8082fbbe6f40.0282d032c43f4d31664d01bc9026ed75a88ebd34.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 18:55:25"

# Widen columns to match the new content.
$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Report for handback generated."
